$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.648550629615784
$ws.Range("B1").Value = 2.980233907699585
$ws.Range("C1").Value = 4.632588386535645
$ws.Range("D1").Value = 1.388772249221802
$ws.Range("E1").Value = 0.8102708458900452
